# Add the "2022-Q1" worksheet (positioned right before the "总计" summary
# sheet), populate it with the quarter's fund-holding detail data, then
# update the "总计" (totals) summary sheet with a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet just before the last sheet
#    ("总计"), matching the tab order in the target workbook.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Use the previous quarter sheet ("2021-Q4") as a formatting template so
# the header row / index column reuse the same bold+border+center style
# already used throughout the workbook.
$template = $wb.Worksheets.Item("2021-Q4")

$template.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Column A (row index) style used throughout the workbook.
$template.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @(0,  "519702", "交银趋势优先混合",            "108.29", "71.40", "2.00", "2.1658", 8),
    @(1,  "671010", "西部利得策略优选混合A",        "6.79",   "93.07", "5.19", "0.3524", 7),
    @(2,  "519979", "长信内需成长混合A",            "6.81",   "92.68", "5.05", "0.3439", 4),
    @(3,  "006397", "长信内需成长混合E",            "6.81",   "92.68", "5.05", "0.3439", 4),
    @(4,  "217024", "招商安盈债券",                "35.05",  "20.20", "0.65", "0.2278", 9),
    @(5,  "011060", "西部利得策略优选混合C",        "1.32",   "93.07", "5.19", "0.0685", 7),
    @(6,  "160620", "鹏华中证A股资源产业指数（LOF）", "2.77",   "94.14", "2.34", "0.0648", 9),
    @(7,  "007423", "西部利得聚禾灵活配置混合A",     "0.60",   "69.21", "4.78", "0.0287", 1),
    @(8,  "007424", "西部利得聚禾灵活配置混合C",     "0.41",   "69.21", "4.78", "0.0196", 1),
    @(9,  "006729", "万家中证500指数增强A",         "1.04",   "93.64", "1.31", "0.0136", 4),
    @(10, "159990", "银华巨潮小盘价值ETF",          "1.06",   "96.39", "1.18", "0.0125", 8),
    @(11, "006730", "万家中证500指数增强C",         "0.61",   "93.64", "1.31", "0.0080", 4),
    @(12, "001657", "长安鑫富领先灵活配置混合",      "0.07",   "30.32", "1.89", "0.0013", 10)
)

# Text-valued columns that look numeric (fund codes / formatted figures)
# need to be force-typed as Text before assignment, otherwise Excel's
# smart type inference would strip leading zeros / trailing zeros and
# store them as numbers. The format is reset back to the default right
# after so no stray number-format style sticks to the cell.
$textCols = @(2, 4, 5, 6, 7)

$r = 2
foreach ($row in $rows) {
    foreach ($col in $textCols) {
        $newSheet.Cells.Item($r, $col).NumberFormat = "@"
    }

    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]

    foreach ($col in $textCols) {
        $newSheet.Cells.Item($r, $col).Style = "Normal"
    }

    $r++
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new summary row for 2022-Q1 above
#    the existing 2021-Q4 row, pushing every other row down by one and
#    renumbering the index column (A) to stay sequential (0..5).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").Style = "Normal"
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 13
$totalSheet.Range("D2").Value = 3.65

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
